$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "25.881.99"
$ws.Range("E2").Value2 = "  -2.27%  "
$ws.Range("D3").Value2 = "1.753.19"
$ws.Range("E3").Value2 = "  -4.64%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  +0.03%  "
$ws.Range("D5").Value = "'239.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -8.26%  "
$ws.Range("E7").Value2 = "  -5.09%  "
$ws.Range("D8").Value = "'42.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -5.62%  "
$ws.Range("D9").Value = "'0.2766"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -5.25%  "
$ws.Range("D10").Value = "'0.06192"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -10.95%  "
$ws.Range("D11").Value2 = "1.749.72"
$ws.Range("E11").Value2 = "  -4.93%  "
$ws.Range("D12").Value = "'15.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -8.43%  "
$ws.Range("D13").Value = "'0.06972"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -3.17%  "
$ws.Range("D14").Value = "'0.6138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -15.21%  "
$ws.Range("D15").Value = "'4.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -8.89%  "
$ws.Range("D16").Value = "'77.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -13.01%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -0.06%  "
$ws.Range("E18").Value2 = "  -0.05%  "
$ws.Range("D19").Value2 = "25.895.55"
$ws.Range("E19").Value2 = "  -2.28%  "
$ws.Range("D20").Value = "'0.000006917"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -12.25%  "
$ws.Range("D21").Value = "'11.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -15.24%  "
$ws.Range("D22").Value2 = "1.972.27"
$ws.Range("E22").Value2 = "  -5.31%  "
$ws.Range("D23").Value = "'4.084"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -10.84%  "
$ws.Range("D24").Value = "'5.262"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -12.06%  "
$ws.Range("D25").Value = "'8.238"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -10.23%  "
$ws.Range("D26").Value = "'138.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -2.58%  "
$ws.Range("D27").Value = "'1.492"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -12.49%  "
$ws.Range("D28").Value = "'1.825"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -15.53%  "
$ws.Range("D29").Value = "'15.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -10.83%  "
$ws.Range("D30").Value = "'103.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -6.28%  "
$ws.Range("D31").Value = "'0.08229"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -7.32%  "
$ws.Range("D32").Value = "'3.697"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -12.97%  "
$ws.Range("D33").Value = "'3.495"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -13.15%  "
$ws.Range("E34").Value2 = "  -6.00%  "
$ws.Range("D35").Value = "'0.9994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -0.07%  "
$ws.Range("D36").Value = "'2.644"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -8.78%  "
$ws.Range("D37").Value = "'0.9923"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -12.17%  "
$ws.Range("D38").Value = "'0.6125"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -15.29%  "
$ws.Range("D39").Value = "'2.709"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -12.33%  "
$ws.Range("D40").Value = "'0.01558"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -8.74%  "
$ws.Range("D41").Value = "'103.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -2.77%  "
$ws.Range("E42").Value2 = "  -0.03%  "
$ws.Range("D43").Value = "'1.892"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -17.67%  "
$ws.Range("D44").Value = "'0.3875"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -16.82%  "
$ws.Range("D45").Value = "'0.7416"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -17.84%  "
$ws.Range("D46").Value = "'4.958"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -15.43%  "
$ws.Range("D47").Value = "'0.05425"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -5.64%  "
$ws.Range("D48").Value = "'0.1115"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -10.12%  "
$ws.Range("D49").Value = "'6.026"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -18.51%  "
$ws.Range("B50").Value2 = "EnergySwap"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.716"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -14.09%  "
$ws.Range("D51").Value = "'52.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -11.57%  "
